$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "nom" column (currently column FZ).
# This shifts "nom" (FZ -> GA) and "url_produit" (GA -> GB) one column to the right,
# matching the dimension change from A1:GA209 to A1:GB209.
$ws.Columns("FZ:FZ").Insert()

# New header cell (row 1) for the freshly inserted price-history column,
# styled the same way as the other header cells (style index 1).
$ws.Range("FZ1").Value = "2026-02-05 10:30:38"

# Price history lookup for rows 2-80: the new column gets a copy of the
# most recent known price (the value that, before the insert, lived in FY
# and now still lives in FY after the shift since FY is to the left of the
# insertion point).
$priceValues = @{
    2 = 39.83
    3 = 169.95
    4 = 169.95
    5 = 249.95
    6 = 299.95
    7 = 339.95
    8 = 619
    9 = 619
    10 = 659
    11 = 659
    12 = 749
    13 = 809
    14 = 809
    15 = 809
    16 = 809
    17 = 809
    18 = 849
    19 = 899
    20 = 899
    21 = 909
    22 = 909
    23 = 909
    24 = 909
    25 = 969
    26 = 969
    27 = 969
    28 = 969
    29 = 969
    30 = 999
    31 = 999
    32 = 1039
    33 = 1039
    34 = 1079
    35 = 1079
    36 = 1079
    37 = 1079
    38 = 1099
    39 = 1099
    40 = 1199
    41 = 1219
    42 = 1219
    43 = 1219
    44 = 1219
    45 = 1219
    46 = 1229
    47 = 1229
    48 = 1249
    49 = 1329
    50 = 1329
    51 = 1329
    52 = 1329
    53 = 1329
    54 = 1329
    55 = 1329
    56 = 1349
    57 = 1419
    58 = 1479
    59 = 1479
    60 = 1479
    61 = 1549
    62 = 1579
    63 = 1579
    64 = 1579
    65 = 1579
    66 = 1579
    67 = 1579
    68 = 1579
    69 = 1729
    70 = 1729
    71 = 1729
    72 = 1829
    73 = 1829
    74 = 1829
    75 = 1979
    76 = 1979
    77 = 1979
    78 = 2479
    79 = 2479
    80 = 2479
}

foreach ($row in $priceValues.Keys) {
    $ws.Cells.Item($row, 182).Value = $priceValues[$row]
}

# Rows 81-209 have no recorded price yet (FY was already blank there), so
# the newly inserted cell stays blank for those rows too.
for ($r = 81; $r -le 209; $r++) {
    $ws.Cells.Item($r, 182).Value = ""
}
